# store version 2 of documents
#
# The document originally contains a single paragraph with the text
# "First version" (and a _GoBack bookmark wrapping it). We need to turn
# this into two paragraphs:
#   1. "First version"            (new paragraph, no bookmark)
#   2. "Second version"           (keeps the original bookmark)
#
# We do this by inserting a brand new paragraph mark before the existing
# paragraph (which duplicates the original paragraph's mark position, so
# the bookmark/extra paragraph metadata stays attached to the *second*
# paragraph, matching how Word behaves when you press Enter at the start
# of a paragraph), filling the new first paragraph with "First version",
# and then updating the original (now second) paragraph's text in place
# (so its run keeps its existing formatting / empty <w:rPr/>) to
# "Second version".

$d = $word.ActiveDocument

# Insert a new, empty paragraph immediately before the first paragraph.
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

# Fill the newly created first paragraph with the original text.
$newFirstPara = $d.Paragraphs.Item(1)
$newFirstPara.Range.Text = "First version"

# The original paragraph (with its bookmark) is now the second paragraph;
# update its text in place to "Second version".
$secondPara = $d.Paragraphs.Item(2)
$secondPara.Range.Text = "Second version"
